$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND (replace): [$old]"
        return
    }
    $r.Text = $new
}

function DeleteText($old) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND (delete): [$old]"
        return
    }
    $r.Delete()
}

# --- Title / byline / email ---
ReplaceText "The Enigmatic Realm of Dark Matter" "Exploring Arts & Literature: Two Halves of a Creative Whole"
ReplaceText "Tomas Nalli" "Harleen Chhabra"
ReplaceText "tomasnalli@post" "harleenchhabrawrites@gmail"

# --- Body paragraph, segment 1 ---
ReplaceText "The cosmic tapestry is a symphony of mysteries, with dark matter composing approximately 27% of its enigmatic composition" "Immerse yourselves in a world where imagination blooms, where brushstrokes and inked words weave together a vibrant symphony of self-expression"
ReplaceText " This uncharted territory has captivated the imaginations of scientists, yearning to unravel its secrets and unlock the hidden treasures of the universe" " In this artistic realm, beauty meets ingenuity, transcending the boundaries of language and form"
ReplaceText " Dark matter's gravitational influence shapes the cosmos, bending and distorting light to create gravitational lensing effects" " It's here that we discover the profound impact of arts and literature; not only as cherished cultural artifacts but also as mirrors reflecting the tapestry of human history and emotion"
ReplaceText " It is the architect of cosmic architecture, responsible for the formation and evolution of galaxies and galaxy clusters" " Unlock the doors to these parallel worlds of creativity, where artists paint with colors, dancers move like liquid poetry, musicians orchestrate emotions, and authors weave narratives that touch the hearts and minds of generations"
ReplaceText " Despite its profound impact, dark matter remains an enigmatic enigma, challenging our understanding of the universe and beckoning us to explore its depths" " Together, let's embark on an exploration of arts and literature, appreciating their profound connections and savoring the unique gifts they bestow upon humanity"

# --- Body paragraph, segment 2 ---
ReplaceText "As we peer into the night sky, the luminous tapestry of stars, galaxies, and nebulae captivates our senses" "From ancient cave paintings to breathtaking sculptures, from elaborate tapestries to modern street art, throughout time, the visual arts have served as a visual language that transcends cultural and linguistic barriers"
ReplaceText " However, there lies a hidden realm, an unseen force that governs the cosmos - dark matter" " Through their strokes and hues, artists chronicle civilizations, express emotions, and explore concepts that defy verbal expression"
ReplaceText " This mysterious substance, comprising approximately 27% of the universe, is invisible to our eyes and instrumentation, revealing its presence only through its gravitational effects" " Whether it be the haunting eyes of the Mona Lisa or the bold lines of a Jackson Pollock, every masterpiece carries a distinct story, inviting viewers to contemplate, interpret, and engage with the artist's inner world"
DeleteText ". The quest to comprehend dark matter has become a scientific odyssey, driving us to push the boundaries of our understanding and unveil the secrets of the universe's composition"

# --- Body paragraph, segment 3 ---
ReplaceText "The influence of dark matter is pervasive, shaping the structure and dynamics of galaxies" "Like the visual arts, literature transports us to distant lands and far-off times, introducing us to myriad characters, cultures, and experiences"
ReplaceText " It governs the motion of stars within galaxies, creating intricate patterns and mesmerizing spiral arms" " With the turn of each page, authors paint pictures with words, using language as their palette"
ReplaceText " Its gravitational pull shapes the destinies of celestial bodies, dictating their orbits and shaping their interactions" " They wield their pens with the dexterity of master craftsmen, forging sentences that evoke emotions, create imagery, and challenge our perceptions of the world"
ReplaceText " Understanding dark matter is not merely an academic pursuit; it holds the key to unlocking the mysteries of the cosmos, revealing the true nature of gravity and the fundamental forces that govern our universe" " From Homer's epics to Shakespeare's sonnets, the written word has the power to ignite imaginations, broaden perspectives, and spark social and intellectual change"

# --- Summary paragraph ---
ReplaceText "The exploration of dark matter is an ongoing scientific journey, pushing the boundaries of our knowledge and understanding of the cosmos" "In conclusion, arts and literature are not merely adornments; they are windows into the human experiences"
ReplaceText " As we delve deeper into the enigmas of this invisible force, we uncover its profound influence on the structure, dynamics, and evolution of the universe" " Together, they provide essential nourishment for our souls, encouraging critical thinking, empathy, and self-awareness"
ReplaceText " The quest to comprehend dark matter will undoubtedly redefine our " " As students of art and literature, we have the privilege of appreciating the splendor of artworks, delving into the depth of literary works, and engaging in discussions that enrich our appreciation of creativity"
ReplaceText "comprehension of gravity, space, and time, painting a more vivid picture of the universe's enigmatic tapestry" "."

# Add the new trailing sentence that used to be the "comprehension of gravity..." run.
$r = $d.Content
$found = $r.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Move the lastRenderedPageBreak from the (old) "comprehension..." run to the new
# first run of the Summary paragraph, and insert the new closing sentence at the
# end of the Summary paragraph (after the "." that used to carry the page break).
$summaryPara = $d.Paragraphs(6)
$summaryRange = $summaryPara.Range
$lastPeriod = $summaryRange.Characters($summaryRange.Characters.Count)
Write-Output "last char: [$($lastPeriod.Text)]"
